$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Enter the text via a formula so Excel's autodetection never turns a
    # date-looking string (e.g. "10/10/2002") into a real date value, then
    # convert the formula result to a plain static value in place so the
    # cell keeps its original (unstyled) text representation.
    $r = $ws.Range($range)
    $escaped = $text.Replace('"', '""')
    $r.Formula = '="' + $escaped + '"'
    $r.Copy() | Out-Null
    $r.PasteSpecial(-4163) | Out-Null
}

# F2: 10/26/02 -> 10/10/2002
Set-TextValue "F2" "10/10/2002"

# E3: martingaido00@gmail.com -> martingaido0@gmail.com
$ws.Range("E3").Value = "martingaido0@gmail.com"

# F3: 19/11/03 -> 18/11/2003
$ws.Range("F3").Value = "18/11/2003"

# F4: 11/10/02 -> 11/10/2002
Set-TextValue "F4" "11/10/2002"

# F5: 2/10/03 -> 2/10/2003
Set-TextValue "F5" "2/10/2003"

# H5: 0 (style 3) -> 100 (style 2)
Set-TextValue "H5" "100"
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# F6: 5/10/89 -> 5/10/1989
Set-TextValue "F6" "5/10/1989"

# F7: 24/10/03 -> 24/10/2003
$ws.Range("F7").Value = "24/10/2003"

# F8: 9/5/77 -> 9/5/1977
Set-TextValue "F8" "9/5/1977"
